$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.488.81"
$ws.Cells.Item(2, 5).Value = "  +5.04%  "
$ws.Cells.Item(3, 4).Value = "1.723.55"
$ws.Cells.Item(3, 5).Value = "  +4.08%  "
$ws.Cells.Item(4, 4).Value = "1.005"
$ws.Cells.Item(4, 5).Value = "  +0.16%  "
$ws.Cells.Item(5, 4).Value = "225.83"
$ws.Cells.Item(5, 5).Value = "  +3.29%  "
$ws.Cells.Item(6, 4).Value = "0.5376"
$ws.Cells.Item(6, 5).Value = "  +2.70%  "
$ws.Cells.Item(8, 4).Value = "0.2682"
$ws.Cells.Item(8, 5).Value = "  +0.69%  "
$ws.Cells.Item(9, 4).Value = "0.06608"
$ws.Cells.Item(9, 5).Value = "  +4.15%  "
$ws.Cells.Item(10, 4).Value = "21.64"
$ws.Cells.Item(10, 5).Value = "  +5.36%  "
$ws.Cells.Item(11, 4).Value = "0.07754"
$ws.Cells.Item(11, 5).Value = "  +0.97%  "
$ws.Cells.Item(12, 4).Value = "4.634"
$ws.Cells.Item(12, 5).Value = "  +0.18%  "
$ws.Cells.Item(13, 4).Value = "1.725.32"
$ws.Cells.Item(13, 5).Value = "  +4.23%  "
$ws.Cells.Item(14, 4).Value = "1.960.40"
$ws.Cells.Item(14, 5).Value = "  +4.00%  "
$ws.Cells.Item(15, 4).Value = "0.5876"
$ws.Cells.Item(15, 5).Value = "  +4.72%  "
$ws.Cells.Item(16, 4).Value = "0.0₅8273"
$ws.Cells.Item(16, 5).Value = "  +1.09%  "
$ws.Cells.Item(17, 4).Value = "67.99"
$ws.Cells.Item(17, 5).Value = "  +3.75%  "
$ws.Cells.Item(18, 4).Value = "27.520.83"
$ws.Cells.Item(18, 5).Value = "  +5.21%  "
$ws.Cells.Item(19, 4).Value = "223.64"
$ws.Cells.Item(19, 5).Value = "  +16.14%  "
$ws.Cells.Item(20, 5).Value = "  +0.06%  "
$ws.Cells.Item(21, 4).Value = "4.736"
$ws.Cells.Item(21, 5).Value = "  +1.62%  "
$ws.Cells.Item(22, 4).Value = "10.71"
$ws.Cells.Item(22, 5).Value = "  +2.20%  "
$ws.Cells.Item(23, 4).Value = "6.102"
$ws.Cells.Item(23, 5).Value = "  +2.40%  "
$ws.Cells.Item(24, 5).Value = "  +0.11%  "
$ws.Cells.Item(25, 4).Value = "148.22"
$ws.Cells.Item(25, 5).Value = "  +1.90%  "
$ws.Cells.Item(26, 5).Value = "  +3.24%  "
$ws.Cells.Item(27, 4).Value = "1.687"
$ws.Cells.Item(27, 5).Value = "  +11.12%  "
$ws.Cells.Item(28, 4).Value = "7.411"
$ws.Cells.Item(28, 5).Value = "  +1.99%  "
$ws.Cells.Item(29, 4).Value = "16.70"
$ws.Cells.Item(29, 5).Value = "  +4.77%  "
$ws.Cells.Item(30, 4).Value = "0.05572"
$ws.Cells.Item(30, 5).Value = "  +1.69%  "
$ws.Cells.Item(31, 4).Value = "1.304"
$ws.Cells.Item(31, 5).Value = "  +2.51%  "
$ws.Cells.Item(32, 5).Value = "  +2.86%  "
$ws.Cells.Item(33, 4).Value = "3.464"
$ws.Cells.Item(33, 5).Value = "  +2.84%  "
$ws.Cells.Item(34, 4).Value = "1.664"
$ws.Cells.Item(34, 5).Value = "  +6.31%  "
$ws.Cells.Item(35, 4).Value = "2.451"
$ws.Cells.Item(35, 5).Value = "  +2.00%  "
$ws.Cells.Item(36, 4).Value = "0.9606"
$ws.Cells.Item(36, 5).Value = "  +1.07%  "
$ws.Cells.Item(37, 4).Value = "2.818"
$ws.Cells.Item(37, 5).Value = "  +1.36%  "
$ws.Cells.Item(38, 5).Value = "  +4.37%  "
$ws.Cells.Item(39, 4).Value = "0.01646"
$ws.Cells.Item(39, 5).Value = "  +3.50%  "
$ws.Cells.Item(40, 4).Value = "5.875"
$ws.Cells.Item(40, 5).Value = "  +0.03%  "
$ws.Cells.Item(41, 4).Value = "0.8569"
$ws.Cells.Item(41, 5).Value = "  +2.76%  "
$ws.Cells.Item(42, 4).Value = "1.058.65"
$ws.Cells.Item(42, 5).Value = "  +2.64%  "
$ws.Cells.Item(43, 4).Value = "1.004"
$ws.Cells.Item(43, 5).Value = "  +0.09%  "
$ws.Cells.Item(44, 4).Value = "101.57"
$ws.Cells.Item(44, 5).Value = "  +0.56%  "
$ws.Cells.Item(45, 5).Value = "  +3.95%  "
$ws.Cells.Item(46, 5).Value = "  +12.02%  "
$ws.Cells.Item(47, 4).Value = "59.02"
$ws.Cells.Item(47, 5).Value = "  +1.55%  "
$ws.Cells.Item(48, 4).Value = "8.222"
$ws.Cells.Item(48, 5).Value = "  +2.55%  "
$ws.Cells.Item(49, 4).Value = "0.4437"
$ws.Cells.Item(49, 5).Value = "  +2.09%  "
$ws.Cells.Item(50, 4).Value = "1.002"
$ws.Cells.Item(50, 5).Value = "  +0.25%  "
$ws.Cells.Item(51, 4).Value = "0.05280"
$ws.Cells.Item(51, 5).Value = "  +1.25%  "